$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.974.71'
$ws.Range('E2').Value = '  +4.05%  '

$ws.Range('D3').Value = '2.781.14'
$ws.Range('E3').Value = '  +4.59%  '

$ws.Range('E4').Value = '  +0.08%  '

$ws.Range('D5').Value = '583.07'
$ws.Range('E5').Value = '  +0.25%  '

$ws.Range('D6').Value = '161.63'
$ws.Range('E6').Value = '  +11.70%  '

$ws.Range('D7').Value = '0.624'
$ws.Range('E7').Value = '  +3.66%  '

$ws.Range('D8').Value = '0.996'
$ws.Range('E8').Value = '  -0.14%  '

$ws.Range('D9').Value = '2.812.18'
$ws.Range('E9').Value = '  +5.00%  '

$ws.Range('D10').Value = '6.86'
$ws.Range('E10').Value = '  +3.74%  '

$ws.Range('D11').Value = '0.115'
$ws.Range('E11').Value = '  +3.75%  '

$ws.Range('E12').Value = '  +4.17%  '

$ws.Range('E13').Value = '  +0.93%  '

$ws.Range('D14').Value = '3.271.86'
$ws.Range('E14').Value = '  +4.60%  '

$ws.Range('D15').Value = '27.82'
$ws.Range('E15').Value = '  +6.48%  '

$ws.Range('D16').Value = '63.969.96'
$ws.Range('E16').Value = '  +4.18%  '

$ws.Range('E17').Value = '  +8.90%  '

$ws.Range('D18').Value = '2.791.42'
$ws.Range('E18').Value = '  +4.81%  '

$ws.Range('D19').Value = '12.36'
$ws.Range('E19').Value = '  +5.67%  '

$ws.Range('D20').Value = '5.03'
$ws.Range('E20').Value = '  +4.89%  '

$ws.Range('D21').Value = '368.75'
$ws.Range('E21').Value = '  +3.51%  '

$ws.Range('D22').Value = '7.11'
$ws.Range('E22').Value = '  +2.69%  '

$ws.Range('D23').Value = '0.550'
$ws.Range('E23').Value = '  +4.50%  '

$ws.Range('E24').Value = '  +0.73%  '

$ws.Range('D25').Value = '67.69'
$ws.Range('E25').Value = '  +4.52%  '

$ws.Range('E26').Value = '  +5.82%  '

$ws.Range('D27').Value = '8.72'
$ws.Range('E27').Value = '  +2.48%  '

$ws.Range('D28').Value = '0.0₃0972'
$ws.Range('E28').Value = '  +17.69%  '

$ws.Range('E29').Value = '  +0.37%  '

$ws.Range('E30').Value = '  +1.68%  '

$ws.Range('D31').Value = '7.33'
$ws.Range('E31').Value = '  +5.29%  '

$ws.Range('D32').Value = '1.27'
$ws.Range('E32').Value = '  +11.44%  '

$ws.Range('D33').Value = '173.34'
$ws.Range('E33').Value = '  +2.15%  '

$ws.Range('D34').Value = '20.97'
$ws.Range('E34').Value = '  +3.75%  '

$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').Value = '5.07'
$ws.Range('E35').Value = '  +8.05%  '

$ws.Range('B36').Value = 'USDe'
$ws.Range('C36').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D36').Value = '0.997'
$ws.Range('E36').Value = '  -0.02%  '

$ws.Range('E37').Value = '  +7.90%  '

$ws.Range('E38').Value = '  +7.28%  '

$ws.Range('E39').Value = '  +2.15%  '

$ws.Range('D40').Value = '343.35'
$ws.Range('E40').Value = '  -0.96%  '

$ws.Range('E41').Value = '  +3.20%  '

$ws.Range('D42').Value = '6.23'
$ws.Range('E42').Value = '  +15.31%  '

$ws.Range('D43').Value = '39.79'
$ws.Range('E43').Value = '  +3.15%  '

$ws.Range('D44').Value = '22.92'
$ws.Range('E44').Value = '  +9.71%  '

$ws.Range('D45').Value = '22.77'
$ws.Range('E45').Value = '  +6.68%  '

$ws.Range('D46').Value = '0.0612'
$ws.Range('E46').Value = '  +5.55%  '

$ws.Range('D47').Value = '0.658'
$ws.Range('E47').Value = '  +4.56%  '

$ws.Range('D48').Value = '0.0262'
$ws.Range('E48').Value = '  +3.02%  '

$ws.Range('D49').Value = '138.39'
$ws.Range('E49').Value = '  +1.89%  '

$ws.Range('E50').Value = '  +2.87%  '

$ws.Range('D51').Value = '2.183.24'
$ws.Range('E51').Value = '  +3.80%  '
